$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (pushes the old row 16.."Lymph nodes / Coming soon"
# and everything below it down by one).
$ws.Rows.Item(16).Insert()

# New entry: Liver vasculature / Congestive Hepatopathy / Clip 1-Bmode + Color Doppler
$ws.Cells.Item(16, 1).Value = "Liver vasculature"
$ws.Cells.Item(16, 2).Value = "Congestive Hepatopathy "
$ws.Cells.Item(16, 3).Value = "Clip 1-Bmode + Color Doppler"

$ws.Hyperlinks.Add($ws.Cells.Item(16, 4), "https://youtu.be/sRu_NTopG3Y", "", "", "https://youtu.be/sRu_NTopG3Y") | Out-Null

# Re-apply the same "hyperlink" cell style used by the other link cells
# (Hyperlinks.Add on its own creates a near-duplicate style).
$ws.Cells.Item(16, 4).Style = $ws.Cells.Item(15, 4).Style

# Match the author's final selection in the saved file.
$ws.Range("D16").Select()
